$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'44.034.43"
$ws.Range("E2").Value = '  +0.94%  '

$ws.Range("D3").Value = "'2.265.09"
$ws.Range("E3").Value = '  +1.12%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").Value = "'318.60"
$ws.Range("E5").Value = '  -0.60%  '

$ws.Range("D6").Value = "'102.43"
$ws.Range("E6").Value = '  +0.84%  '

$ws.Range("E7").Value = '  +0.39%  '

$ws.Range("E8").Value = '  -0.06%  '

$ws.Range("D9").Value = "'0.570"
$ws.Range("E9").Value = '  -0.34%  '

$ws.Range("D10").Value = "'38.32"
$ws.Range("E10").Value = '  +3.26%  '

$ws.Range("D11").Value = "'0.0837"
$ws.Range("E11").Value = '  +0.89%  '

$ws.Range("D12").Value = "'7.83"
$ws.Range("E12").Value = '  +0.71%  '

$ws.Range("D13").Value = "'0.107"
$ws.Range("E13").Value = '  +1.37%  '

$ws.Range("D14").Value = "'2.612.59"
$ws.Range("E14").Value = '  +1.53%  '

$ws.Range("D15").Value = "'0.874"
$ws.Range("E15").Value = '  +0.49%  '

$ws.Range("D16").Value = "'14.54"
$ws.Range("E16").Value = '  +2.11%  '

$ws.Range("D17").Value = "'2.271.25"
$ws.Range("E17").Value = '  +1.35%  '

$ws.Range("D18").Value = "'43.916.02"
$ws.Range("E18").Value = '  +1.43%  '

$ws.Range("D19").Value = "'14.41"
$ws.Range("E19").Value = '  -0.24%  '

$ws.Range("D20").Value = "'0.0₃0989"
$ws.Range("E20").Value = '  +0.61%  '

$ws.Range("D21").Value = "'6.64"
$ws.Range("E21").Value = '  +0.94%  '

$ws.Range("D22").Value = "'65.98"
$ws.Range("E22").Value = '  +0.03%  '

$ws.Range("D23").Value = "'3.21"
$ws.Range("E23").Value = '  -0.76%  '

$ws.Range("D24").Value = "'238.49"
$ws.Range("E24").Value = '  +0.21%  '

$ws.Range("D25").Value = "'2.19"
$ws.Range("E25").Value = '  +0.76%  '

$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = '  +0.02%  '

$ws.Range("E27").Value = '  +1.40%  '

$ws.Range("E28").Value = '  -0.69%  '

$ws.Range("D29").Value = "'38.71"
$ws.Range("E29").Value = '  +11.95%  '

$ws.Range("E30").Value = '  +0.13%  '

$ws.Range("D31").Value = "'6.53"
$ws.Range("E31").Value = '  +0.22%  '

$ws.Range("D32").Value = "'163.91"
$ws.Range("E32").Value = '  +3.40%  '

$ws.Range("D33").Value = "'0.0881"
$ws.Range("E33").Value = '  -1.87%  '

$ws.Range("D34").Value = "'20.45"
$ws.Range("E34").Value = '  -1.52%  '

$ws.Range("E35").Value = '  -2.27%  '

$ws.Range("D36").Value = "'3.28"
$ws.Range("E36").Value = '  -4.37%  '

$ws.Range("D37").Value = "'2.02"
$ws.Range("E37").Value = '  +2.10%  '

$ws.Range("E38").Value = '  -1.87%  '

$ws.Range("E39").Value = '  +0.04%  '

$ws.Range("D40").Value = "'0.109"
$ws.Range("E40").Value = '  +2.63%  '

$ws.Range("D41").Value = "'3.85"
$ws.Range("E41").Value = '  +4.75%  '

$ws.Range("D42").Value = "'15.74"
$ws.Range("E42").Value = '  +28.50%  '

$ws.Range("D43").Value = "'0.0325"
$ws.Range("E43").Value = '  -0.76%  '

$ws.Range("E44").Value = '  +0.23%  '

$ws.Range("D45").Value = "'1.776.73"
$ws.Range("E45").Value = '  -2.66%  '

$ws.Range("D46").Value = "'0.207"
$ws.Range("E46").Value = '  -1.55%  '

$ws.Range("D47").Value = "'85.00"
$ws.Range("E47").Value = '  -5.56%  '

$ws.Range("D48").Value = "'5.39"
$ws.Range("E48").Value = '  -3.14%  '

$ws.Range("D49").Value = "'8.94"
$ws.Range("E49").Value = '  +4.28%  '

$ws.Range("B50").Value = 'MultiversX'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D50").Value = "'59.49"
$ws.Range("E50").Value = '  -3.52%  '

$ws.Range("B51").Value = 'ordi'
$ws.Range("C51").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D51").Value = "'74.61"
$ws.Range("E51").Value = '  -7.31%  '
